# Chapter 5 edit: rename three TOC bookmarks and bump two displayed
# page numbers (header1.xml "20" -> "50", footer1.xml "19" -> "49").

$d = $word.ActiveDocument

# --- 1. Rename hidden _Toc bookmarks -----------------------------------
# Word's COM model does not let us rewrite Bookmark.Name in place in this
# runtime, so recreate each bookmark under its new name over the exact
# same range, then delete the old one (net effect identical to a rename,
# and the implementation re-uses the freed bookmark id on save).
function Rename-Bookmark($oldName, $newName) {
    $bm = $d.Bookmarks($oldName)
    $r = $bm.Range
    $d.Bookmarks.Add($newName, $r)
    $d.Bookmarks($oldName).Delete()
}

Rename-Bookmark "_Toc428458295" "_Toc430350704"
Rename-Bookmark "_Toc428458296" "_Toc430350705"
Rename-Bookmark "_Toc428458297" "_Toc430350706"

# --- 2. Update the displayed page numbers in header/footer -------------
# header1.xml / footer1.xml belong to the first section (titlePg: the
# "default" header is index 1, the "first page" footer is index 2).
$sec1 = $d.Sections(1)

$hdr = $sec1.Headers(1)
$hdr.Range.Find.Execute("20", $true, $false, $false, $false, $false, $true, 1, $false, "50", 2)

$ftr = $sec1.Footers(2)
$ftr.Range.Find.Execute("19", $true, $false, $false, $false, $false, $true, 1, $false, "49", 2)
